# Apply cryptos list price/volume updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "95.616.30"
    "E2" = "  -0.53%  "
    "D3" = "3.591.35"
    "E3" = "  -2.14%  "
    "E4" = "  -0.02%  "
    "D5" = "237.52"
    "E5" = "  -1.37%  "
    "D6" = "659.07"
    "E6" = "  +2.28%  "
    "E7" = "  +1.29%  "
    "D8" = "0.404"
    "E8" = "  +0.20%  "
    "E9" = "  +0.07%  "
    "D10" = "1.04"
    "E10" = "  +2.90%  "
    "D11" = "3.591.55"
    "E11" = "  -2.14%  "
    "D12" = "43.03"
    "E12" = "  -1.98%  "
    "E13" = "  +0.79%  "
    "D14" = "6.50"
    "E14" = "  +1.94%  "
    "D15" = "4.260.33"
    "E15" = "  -2.21%  "
    "D16" = "95.485.82"
    "E16" = "  -0.54%  "
    "E17" = "  -0.40%  "
    "D18" = "3.592.20"
    "E18" = "  -2.22%  "
    "D19" = "12.86"
    "E19" = "  -6.15%  "
    "E20" = "  -7.95%  "
    "D21" = "18.03"
    "E21" = "  -3.67%  "
    "D22" = "3.48"
    "E22" = "  +0.87%  "
    "E23" = "  +1.27%  "
    "D24" = "511.84"
    "E24" = "  -1.41%  "
    "D25" = "7.15"
    "E25" = "  +5.09%  "
    "E26" = "  -0.02%  "
    "D27" = "95.98"
    "E27" = "  -1.91%  "
    "D28" = "12.90"
    "E28" = "  +2.51%  "
    "D29" = "3.785.51"
    "E29" = "  -1.96%  "
    "D30" = "3.07"
    "E30" = "  -4.10%  "
    "E31" = "  +3.05%  "
    "D32" = "11.62"
    "E32" = "  -0.76%  "
    "D34" = "1.00"
    "E34" = "  -0.48%  "
    "E35" = "  -1.45%  "
    "D36" = "32.18"
    "E36" = "  -1.04%  "
    "E37" = "  +14.01%  "
    "D38" = "8.73"
    "E38" = "  +10.77%  "
    "E39" = "  -2.43%  "
    "D40" = "602.06"
    "E40" = "  +6.24%  "
    "E41" = "  -0.11%  "
    "E42" = "  +0.13%  "
    "E43" = "  +7.97%  "
    "E44" = "  -4.25%  "
    "D45" = "35.30"
    "E45" = "  +9.14%  "
    "E46" = "  +0.20%  "
    "E47" = "  +3.45%  "
    "E48" = "  -2.76%  "
    "E49" = "  -1.34%  "
    "E50" = "  +0.09%  "
    "D51" = "8.25"
    "E51" = "  -1.01%  "
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    $range.Style = "Normal"
}

